# Generate Report for Handoff
# Mark the feecfe05-1d7e-4053-9c6d-2583a85d1160 file as "Ready for handoff"
# across the Overview, zh-cn and de-de sheets, and refresh its handoff
# datetime stamps.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet - row 3 corresponds to feecfe05-1d7e-4053-9c6d-2583a85d1160.md
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-15-20 04:15:15"

# zh-cn sheet - row 3 corresponds to feecfe05-1d7e-4053-9c6d-2583a85d1160.md
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "2016-03-20 04:15:12"

# de-de sheet - row 3 corresponds to feecfe05-1d7e-4053-9c6d-2583a85d1160.md
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "2016-03-20 04:15:15"
